# The observation rows 34-50 got re-shuffled into a new order (same 17
# records, different row positions). Snapshot every row's contents first,
# then write each row back out at its new position according to the
# mapping (new row -> source/old row).
#
# Columns I, AT and AY hold an explicit-but-empty inline string in every
# one of these rows both before and after the edit, so they are left
# completely untouched (round-tripping a blank through Value/Value2 drops
# the cell instead of keeping the empty placeholder). Likewise the many
# columns that have no cell at all in any of these rows (J-L, N-O, X, AC,
# AF, AH-AS, AU-AV, ...) are skipped so no stray blank cells get created.
# What's left are these contiguous column blocks that actually carry data:

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 34
$lastRow = 50
$colBlocks = @("A:H", "M:M", "P:W", "Z:Z", "AB:AB", "AD:AE", "AG:AG", "AW:AX")
# Y and AA hold plain "yyyy-mm-dd" text (General format, stored as a
# string, not a real date). They need special handling below, because
# writing a date-shaped string straight into Value/Value2 makes Excel
# "helpfully" reinterpret it as a real date serial number.
$dateTextCols = @("Y", "AA")

# new row -> old row it should now contain
$mapping = @{
    34 = 38
    35 = 42
    36 = 34
    37 = 43
    38 = 45
    39 = 39
    40 = 37
    41 = 41
    42 = 35
    43 = 48
    44 = 47
    45 = 49
    46 = 40
    47 = 36
    48 = 44
    49 = 50
    50 = 46
}

# Snapshot every block of every row before overwriting anything, since
# several rows are both sources and destinations.
$snapshots = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowBlocks = @{}
    foreach ($block in $colBlocks) {
        $cols = $block.Split(":")
        $rowBlocks[$block] = $ws.Range($cols[0] + $r + ":" + $cols[1] + $r).Value2
    }
    foreach ($col in $dateTextCols) {
        $rowBlocks[$col] = $ws.Range($col + $r).Value2
    }
    $snapshots[$r] = $rowBlocks
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $rowBlocks = $snapshots[$oldRow]
    foreach ($block in $colBlocks) {
        $cols = $block.Split(":")
        $ws.Range($cols[0] + $newRow + ":" + $cols[1] + $newRow).Value2 = $rowBlocks[$block]
    }
    foreach ($col in $dateTextCols) {
        $target = $ws.Range($col + $newRow)
        # Prefix with an apostrophe to force text, which stops Excel from
        # auto-converting "2023-08-22" into a date serial, then reset the
        # style so no extra quote-prefix/number-format style lingers.
        $target.Value = "'" + $rowBlocks[$col]
        $target.Style = "Normal"
    }
}
